$d = $word.ActiveDocument
$d.Content.Find.Execute("Исполнитель", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Разработчик", 2)
